$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Component T1's Comment (A12) changes from "BS170" to "2N7002"
$ws.Range("A12").Value = "2N7002"

# A7/B7 had a stray "applyFill" style that differs only cosmetically from the
# style used by every other data row; normalise them onto the same format as
# the rest of the table (no explicit fill application).
$ws.Range("A7:B7").Interior.Pattern = -4142

# The active selection at save time was A13.
$ws.Range("A13").Select() | Out-Null
